$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cryptocurrency price/volume updates scraped on Mon Mar 27 05:36:34 UTC 2023.
# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (plain decimals like "1.003") are forced to Text format first so the exact
# scraped string (with its original precision/formatting) is preserved verbatim,
# matching the inline-string cell type used throughout this sheet.

$ws.Cells.Item(2, 4).Value = "27.907.09"  # D2: '27.923.91' -> '27.907.09'
$ws.Cells.Item(2, 5).Value = "  +1.28%  "  # E2: '  +1.32%  ' -> '  +1.28%  '
$ws.Cells.Item(3, 4).Value = "1.765.06"  # D3: '1.764.97' -> '1.765.06'
$ws.Cells.Item(3, 5).Value = "  +0.77%  "  # E3: '  +0.78%  ' -> '  +0.77%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"  # D4: force text so "1.003" is not parsed as a number
$ws.Cells.Item(4, 4).Value = "1.003"  # D4: '1.002' -> '1.003'
$ws.Cells.Item(4, 5).Value = "  +0.00%  "  # E4: '  -0.04%  ' -> '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"  # D5: force text so "329.03" is not parsed as a number
$ws.Cells.Item(5, 4).Value = "329.03"  # D5: '329.12' -> '329.03'
$ws.Cells.Item(5, 5).Value = "  +1.50%  "  # E5: '  +1.53%  ' -> '  +1.50%  '
$ws.Cells.Item(6, 5).Value = "  -0.03%  "  # E6: '  -0.04%  ' -> '  -0.03%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"  # D7: force text so "0.4541" is not parsed as a number
$ws.Cells.Item(7, 4).Value = "0.4541"  # D7: '0.4544' -> '0.4541'
$ws.Cells.Item(7, 5).Value = "  +1.10%  "  # E7: '  +1.25%  ' -> '  +1.10%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"  # D8: force text so "0.3511" is not parsed as a number
$ws.Cells.Item(8, 4).Value = "0.3511"  # D8: '0.3516' -> '0.3511'
$ws.Cells.Item(8, 5).Value = "  -1.36%  "  # E8: '  -1.24%  ' -> '  -1.36%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"  # D9: force text so "42.03" is not parsed as a number
$ws.Cells.Item(9, 4).Value = "42.03"  # D9: '42.05' -> '42.03'
$ws.Cells.Item(10, 4).NumberFormat = "@"  # D10: force text so "0.07384" is not parsed as a number
$ws.Cells.Item(10, 4).Value = "0.07384"  # D10: '0.07382' -> '0.07384'
$ws.Cells.Item(10, 5).Value = "  -0.92%  "  # E10: '  -0.96%  ' -> '  -0.92%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"  # D11: force text so "1.096" is not parsed as a number
$ws.Cells.Item(11, 4).Value = "1.096"  # D11: '1.097' -> '1.096'
$ws.Cells.Item(11, 5).Value = "  +1.40%  "  # E11: '  +1.44%  ' -> '  +1.40%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"  # D12: force text so "1.001" is not parsed as a number
$ws.Cells.Item(12, 4).Value = "1.001"  # D12: '1.002' -> '1.001'
$ws.Cells.Item(12, 5).Value = "  -0.02%  "  # E12: '  +0.02%  ' -> '  -0.02%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"  # D13: force text so "20.70" is not parsed as a number
$ws.Cells.Item(13, 4).Value = "20.70"  # D13: '20.72' -> '20.70'
$ws.Cells.Item(13, 5).Value = "  -0.06%  "  # E13: '  +0.01%  ' -> '  -0.06%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"  # D14: force text so "5.998" is not parsed as a number
$ws.Cells.Item(14, 4).Value = "5.998"  # D14: '5.995' -> '5.998'
$ws.Cells.Item(14, 5).Value = "  +0.38%  "  # E14: '  +0.41%  ' -> '  +0.38%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"  # D15: force text so "7.186" is not parsed as a number
$ws.Cells.Item(15, 4).Value = "7.186"  # D15: '7.187' -> '7.186'
$ws.Cells.Item(15, 5).Value = "  +0.60%  "  # E15: '  +0.63%  ' -> '  +0.60%  '
$ws.Cells.Item(16, 4).Value = "1.770.45"  # D16: '1.767.60' -> '1.770.45'
$ws.Cells.Item(16, 5).Value = "  +1.02%  "  # E16: '  +0.94%  ' -> '  +1.02%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"  # D17: force text so "92.36" is not parsed as a number
$ws.Cells.Item(17, 4).Value = "92.36"  # D17: '92.38' -> '92.36'
$ws.Cells.Item(17, 5).Value = "  -1.68%  "  # E17: '  -1.34%  ' -> '  -1.68%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"  # D18: force text so "0.00001057" is not parsed as a number
$ws.Cells.Item(18, 4).Value = "0.00001057"  # D18: '0.00001058' -> '0.00001057'
$ws.Cells.Item(18, 5).Value = "  +0.26%  "  # E18: '  +0.34%  ' -> '  +0.26%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"  # D19: force text so "0.06448" is not parsed as a number
$ws.Cells.Item(19, 4).Value = "0.06448"  # D19: '0.06451' -> '0.06448'
$ws.Cells.Item(19, 5).Value = "  +1.09%  "  # E19: '  +1.35%  ' -> '  +1.09%  '
$ws.Cells.Item(20, 5).Value = "  +0.00%  "  # E20: '  +0.01%  ' -> '  +0.00%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"  # D21: force text so "16.96" is not parsed as a number
$ws.Cells.Item(21, 4).Value = "16.96"  # D21: '16.95' -> '16.96'
$ws.Cells.Item(21, 5).Value = "  -0.31%  "  # E21: '  -0.61%  ' -> '  -0.31%  '
$ws.Cells.Item(22, 5).Value = "  +0.87%  "  # E22: '  +0.80%  ' -> '  +0.87%  '
$ws.Cells.Item(23, 4).Value = "27.937.31"  # D23: '27.949.35' -> '27.937.31'
$ws.Cells.Item(23, 5).Value = "  +1.22%  "  # E23: '  +1.24%  ' -> '  +1.22%  '
$ws.Cells.Item(24, 5).Value = "  +0.51%  "  # E24: '  +0.45%  ' -> '  +0.51%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"  # D25: force text so "2.154" is not parsed as a number
$ws.Cells.Item(25, 4).Value = "2.154"  # D25: '2.151' -> '2.154'
$ws.Cells.Item(25, 5).Value = "  +3.37%  "  # E25: '  +3.16%  ' -> '  +3.37%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"  # D26: force text so "162.09" is not parsed as a number
$ws.Cells.Item(26, 4).Value = "162.09"  # D26: '162.32' -> '162.09'
$ws.Cells.Item(26, 5).Value = "  -1.95%  "  # E26: '  -1.89%  ' -> '  -1.95%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"  # D27: force text so "20.14" is not parsed as a number
$ws.Cells.Item(27, 4).Value = "20.14"  # D27: '20.15' -> '20.14'
$ws.Cells.Item(27, 5).Value = "  +0.13%  "  # E27: '  +0.16%  ' -> '  +0.13%  '
$ws.Cells.Item(28, 4).Value = "1.971.90"  # D28: '1.972.61' -> '1.971.90'
$ws.Cells.Item(28, 5).Value = "  +0.96%  "  # E28: '  +0.98%  ' -> '  +0.96%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"  # D29: force text so "2.170" is not parsed as a number
$ws.Cells.Item(29, 4).Value = "2.170"  # D29: '2.157' -> '2.170'
$ws.Cells.Item(29, 5).Value = "  +3.42%  "  # E29: '  +2.66%  ' -> '  +3.42%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"  # D30: force text so "123.85" is not parsed as a number
$ws.Cells.Item(30, 4).Value = "123.85"  # D30: '123.84' -> '123.85'
$ws.Cells.Item(30, 5).Value = "  -1.07%  "  # E30: '  -1.03%  ' -> '  -1.07%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"  # D31: force text so "1.073" is not parsed as a number
$ws.Cells.Item(31, 4).Value = "1.073"  # D31: '1.072' -> '1.073'
$ws.Cells.Item(31, 5).Value = "  -0.81%  "  # E31: '  -1.05%  ' -> '  -0.81%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"  # D32: force text so "0.09285" is not parsed as a number
$ws.Cells.Item(32, 4).Value = "0.09285"  # D32: '0.09296' -> '0.09285'
$ws.Cells.Item(32, 5).Value = "  +1.40%  "  # E32: '  +1.41%  ' -> '  +1.40%  '
$ws.Cells.Item(33, 2).Value = "Filecoin"  # B33: 'HuobiToken' -> 'Filecoin'
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"  # C33: 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' -> 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(33, 4).NumberFormat = "@"  # D33: force text so "5.583" is not parsed as a number
$ws.Cells.Item(33, 4).Value = "5.583"  # D33: '3.667' -> '5.583'
$ws.Cells.Item(33, 5).Value = "  +1.84%  "  # E33: '  +0.39%  ' -> '  +1.84%  '
$ws.Cells.Item(34, 2).Value = "HuobiToken"  # B34: 'Filecoin' -> 'HuobiToken'
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"  # C34: 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' -> 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(34, 4).NumberFormat = "@"  # D34: force text so "3.654" is not parsed as a number
$ws.Cells.Item(34, 4).Value = "3.654"  # D34: '5.583' -> '3.654'
$ws.Cells.Item(34, 5).Value = "  +0.05%  "  # E34: '  +1.70%  ' -> '  +0.05%  '
$ws.Cells.Item(35, 5).Value = "  +1.17%  "  # E35: '  +0.86%  ' -> '  +1.17%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"  # D36: force text so "0.02274" is not parsed as a number
$ws.Cells.Item(36, 4).Value = "0.02274"  # D36: '0.02275' -> '0.02274'
$ws.Cells.Item(36, 5).Value = "  -0.32%  "  # E36: '  -0.27%  ' -> '  -0.32%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"  # D37: force text so "0.06116" is not parsed as a number
$ws.Cells.Item(37, 4).Value = "0.06116"  # D37: '0.06121' -> '0.06116'
$ws.Cells.Item(37, 5).Value = "  +2.02%  "  # E37: '  +1.92%  ' -> '  +2.02%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"  # D38: force text so "0.2086" is not parsed as a number
$ws.Cells.Item(38, 4).Value = "0.2086"  # D38: '0.2085' -> '0.2086'
$ws.Cells.Item(38, 5).Value = "  -0.06%  "  # E38: '  -0.02%  ' -> '  -0.06%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"  # D39: force text so "4.933" is not parsed as a number
$ws.Cells.Item(39, 4).Value = "4.933"  # D39: '4.935' -> '4.933'
$ws.Cells.Item(39, 5).Value = "  +0.56%  "  # E39: '  +0.33%  ' -> '  +0.56%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"  # D40: force text so "0.6251" is not parsed as a number
$ws.Cells.Item(40, 4).Value = "0.6251"  # D40: '0.6249' -> '0.6251'
$ws.Cells.Item(40, 5).Value = "  -0.18%  "  # E40: '  -0.28%  ' -> '  -0.18%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"  # D41: force text so "1.181" is not parsed as a number
$ws.Cells.Item(41, 4).Value = "1.181"  # D41: '1.183' -> '1.181'
$ws.Cells.Item(41, 5).Value = "  +0.24%  "  # E41: '  +0.54%  ' -> '  +0.24%  '
$ws.Cells.Item(42, 5).Value = "  -0.95%  "  # E42: '  -1.05%  ' -> '  -0.95%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"  # D43: force text so "7.865" is not parsed as a number
$ws.Cells.Item(43, 4).Value = "7.865"  # D43: '7.807' -> '7.865'
$ws.Cells.Item(43, 5).Value = "  +1.30%  "  # E43: '  +0.56%  ' -> '  +1.30%  '
$ws.Cells.Item(44, 5).Value = "  -0.28%  "  # E44: '  -0.29%  ' -> '  -0.28%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"  # D45: force text so "3.733" is not parsed as a number
$ws.Cells.Item(45, 4).Value = "3.733"  # D45: '3.735' -> '3.733'
$ws.Cells.Item(45, 5).Value = "  +0.50%  "  # E45: '  +0.54%  ' -> '  +0.50%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"  # D46: force text so "0.5838" is not parsed as a number
$ws.Cells.Item(46, 4).Value = "0.5838"  # D46: '0.5845' -> '0.5838'
$ws.Cells.Item(46, 5).Value = "  -0.04%  "  # E46: '  -0.11%  ' -> '  -0.04%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"  # D47: force text so "122.68" is not parsed as a number
$ws.Cells.Item(47, 4).Value = "122.68"  # D47: '122.64' -> '122.68'
$ws.Cells.Item(47, 5).Value = "  +0.83%  "  # E47: '  +0.72%  ' -> '  +0.83%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"  # D48: force text so "1.934" is not parsed as a number
$ws.Cells.Item(48, 4).Value = "1.934"  # D48: '1.935' -> '1.934'
$ws.Cells.Item(48, 5).Value = "  +0.32%  "  # E48: '  +0.24%  ' -> '  +0.32%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"  # D49: force text so "1.127" is not parsed as a number
$ws.Cells.Item(49, 4).Value = "1.127"  # D49: '1.128' -> '1.127'
$ws.Cells.Item(49, 5).Value = "  +0.14%  "  # E49: '  +0.13%  ' -> '  +0.14%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"  # D50: force text so "0.06807" is not parsed as a number
$ws.Cells.Item(50, 4).Value = "0.06807"  # D50: '0.06810' -> '0.06807'
$ws.Cells.Item(50, 5).Value = "  -1.04%  "  # E50: '  -1.02%  ' -> '  -1.04%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"  # D51: force text so "72.92" is not parsed as a number
$ws.Cells.Item(51, 4).Value = "72.92"  # D51: '72.85' -> '72.92'
$ws.Cells.Item(51, 5).Value = "  +2.28%  "  # E51: '  +2.14%  ' -> '  +2.28%  '
